# Updated cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.780.40"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "1.635.88"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  -0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0639"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "1.860.57"
$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("D14").Value = "1.633.94"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.561"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").Value = "25.806.26"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("E25").Value = "  +3.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.18%  "

$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  +1.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("E32").Value = "  +1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("D37").Value = "1.132.07"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "

$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0157"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("E42").Value = "  +1.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.807"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("D45").Value = "1.770.72"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.74%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.95%  "
